# Station.pptx edit: flatten the two "folded corner" callout shapes on
# slide 2 (the Sub/StationLH and Sub/StationRH labels) back to a
# square corner by resetting their adjustment handle to 0.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# "사각형: 모서리가 접힌 도형 2" (folded-corner rectangle 2)
$shp1 = $s.Shapes.Item(5)
$shp1.Adjustments.Item(1) = 0

# "사각형: 모서리가 접힌 도형 4" (folded-corner rectangle 4)
$shp2 = $s.Shapes.Item(6)
$shp2.Adjustments.Item(1) = 0
